$wb = $excel.ActiveWorkbook

$wsRankings = $wb.Worksheets.Item("rankings")
$wsHosts = $wb.Worksheets.Item("hosts")

# --- rankings sheet: add row 27 for year 2023 ---
$wsRankings.Range("A27").Value = 2023
$wsRankings.Range("B27").Value = 14
$wsRankings.Range("D27").Value = 1
$wsRankings.Range("E27").Value = 5
$wsRankings.Range("F27").Value = 8
$wsRankings.Range("G27").Value = 2
$wsRankings.Range("H27").Value = 10
$wsRankings.Range("I27").Value = 7
$wsRankings.Range("J27").Value = 12
$wsRankings.Range("L27").Value = 15
$wsRankings.Range("O27").Value = 11
$wsRankings.Range("Q27").Value = 3
$wsRankings.Range("R27").Value = 13
$wsRankings.Range("U27").Value = 9
$wsRankings.Range("V27").Value = 16
$wsRankings.Range("W27").Value = 6
$wsRankings.Range("Y27").Value = 4

# --- hosts sheet: add row 26 for year 2022, row 27 for year 2023 ---
$wsHosts.Range("A26").Value = 2022
$wsHosts.Range("B26").Value = 0
$wsHosts.Range("C26").Value = 0
$wsHosts.Range("D26").Value = 0
$wsHosts.Range("E26").Value = 0
$wsHosts.Range("F26").Value = 0
$wsHosts.Range("G26").Value = 0
$wsHosts.Range("H26").Value = 0
$wsHosts.Range("I26").Value = 1
$wsHosts.Range("J26").Value = 0
$wsHosts.Range("K26").Value = 0
$wsHosts.Range("L26").Value = 0
$wsHosts.Range("M26").Value = 0
$wsHosts.Range("N26").Value = 0
$wsHosts.Range("O26").Value = 0
$wsHosts.Range("P26").Value = 0
$wsHosts.Range("Q26").Value = 0
$wsHosts.Range("R26").Value = 0
$wsHosts.Range("S26").Value = 0
$wsHosts.Range("T26").Value = 0
$wsHosts.Range("U26").Value = 0
$wsHosts.Range("V26").Value = 0
$wsHosts.Range("W26").Value = 0
$wsHosts.Range("X26").Value = 0
$wsHosts.Range("Y26").Value = 0

$wsHosts.Range("A27").Value = 2023
$wsHosts.Range("B27").Value = 0
$wsHosts.Range("C27").Value = 0
$wsHosts.Range("D27").Value = 0
$wsHosts.Range("E27").Value = 0
$wsHosts.Range("F27").Value = 0
$wsHosts.Range("G27").Value = 0
$wsHosts.Range("H27").Value = 0
$wsHosts.Range("I27").Value = 1
$wsHosts.Range("J27").Value = 0
$wsHosts.Range("K27").Value = 0
$wsHosts.Range("L27").Value = 0
$wsHosts.Range("M27").Value = 0
$wsHosts.Range("N27").Value = 0
$wsHosts.Range("O27").Value = 0
$wsHosts.Range("P27").Value = 0
$wsHosts.Range("Q27").Value = 1
$wsHosts.Range("R27").Value = 0
$wsHosts.Range("S27").Value = 0
$wsHosts.Range("T27").Value = 0
$wsHosts.Range("U27").Value = 0
$wsHosts.Range("V27").Value = 0
$wsHosts.Range("W27").Value = 0
$wsHosts.Range("X27").Value = 0
$wsHosts.Range("Y27").Value = 0

# --- selection / active sheet updates ---
# rankings: cursor moved to A47 but sheet is no longer the active tab
[void]$wsRankings.Range("A47").Select()
# hosts becomes the active/selected sheet, cursor at C41
[void]$wsHosts.Select()
[void]$wsHosts.Range("C41").Select()
